$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.870.94'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").Value = '1.809.81'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.50'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4635'
$ws.Range("E7").Value = '  +3.77%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3703'
$ws.Range("E8").Value = '  -1.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07355'
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8767'
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.47'
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").Value = '1.792.85'
$ws.Range("E12").Value = '  -1.94%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.364'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.523'
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.73'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07043'
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008695'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.73'
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").Value = '26.887.90'
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.320'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("E23").Value = '  -3.19%  '
$ws.Range("D24").Value = '2.038.92'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.899'
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.61'
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.40'
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.156'
$ws.Range("E28").Value = '  -4.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.335'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.23'
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08911'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7528'
$ws.Range("E32").Value = '  -4.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.155'
$ws.Range("E33").Value = '  -3.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.467'
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.923'
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.102'
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01966'
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.424'
$ws.Range("E40").Value = '  +3.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.925'
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.170'
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1666'
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.488'
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4979'
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.31'
$ws.Range("E47").Value = '  -3.07%  '
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.671'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.67'
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06289'
$ws.Range("E51").Value = '  -1.42%  '
